$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headline = @'
boring stock offering an over annual dividend yield there a small problem with trump export deal with nvidia and amd the constitution say it illegal fortune trump chief intervened to save rfk jr top vaccine aide politico revisiting carrie bradshaws apartment in sex and the city and and just like that fed face high bar for big cut despite white house pressure the new york time fed face high bar for big cut despite white house pressure bitcoin price soar to a etf inflow surge how trump get what he want from company nbc news wanf in atlanta release new schedule for life after cbs home for sale in manhattan and brooklyn the live music business is booming now rap is getting a piece too there money to be made in maha food company want in cetera investment adviser sell share of camden property trust nysecpt the cost to rebuild a home after a natural disaster is rising pnc financial service group inc acquires share of on holding ag nyseonon pnc financial service group inc raise position in five below inc nasdaqfive pnc financial service group inc purchase share of the madison square garden company nysemsgs tko group holding inc nysetko share acquired by pnc financial service group inc the latest chatgpt is supposed to be phd level smart it cant even label a map crypto presale alert blockchainfx bfx set for explosive growth can it beat bitcoin hyper blockdag big tech ai boom is reordering the u power grid the new york time crisis intervention effort yield impressive result trump is tightening the screw on corporate america and ceo are staying mum titanbay appoints exubs exec a head of switzerland is xrp a buy while it below best spy tv show to watch now on peacock in vipshop report unaudited second quarter financial result instacart is partying like it with third eye blind frosted tip and ring pop new report cleaner manufacturing could save life and prevent million of asthma attack titanbay appoints exubs exec a head of switzerland trump is tightening the screw on corporate america and ceo are staying mum npr weibo announces second quarter unaudited financial result ai development is accelerating vanguard etf to buy right now inside thailand frantic effort to close tariff loophole to the u china offer interest subsidy for loan amid sluggish economy pinterest share just sank time to buy the dip or run for the hill cava group nysecava price target cut to by analyst at morgan stanley nvidia corporation nasdaqnvda share purchased by richard c young co ltd kingsview wealth management llc purchase share of nvidia corporation nasdaqnvda cava group nysecava price target lowered to at citigroup tat technology tatt buy rating reaffirmed at lake street capital summit x llc purchase share of amazoncom inc nasdaqamzn nvidia corporation nasdaqnvda share sold by forum financial management lp repay nasdaqrpay receives neutral rating from ubs group snowflake nysesnow given new price target at ubs group nvidia corporation nasdaqnvda share sold by asset advisor investment management llc rivian automotive nasdaqrivn given new price target at morgan stanley fluence energy nasdaqflnc price target lowered to at royal bank of canada commerzbank aktiengesellschaft fi ha million stock holding in nvidia corporation nasdaqnvda meet the yield dividend stock that could soar in brachioplasty market to surpass usd billion by fueled by rising aesthetic awareness and technological advancement sn insider titanbay appoints paul essen bakery owner say they closed for a mental health break their employee say theyve been ghosted wkmg in orlando give news director autumn jones all the job vanguard group inc trim stock holding in blueprint medicine corporation nasdaqbpmc vanguard group inc ha million stock position in interdigital inc nasdaqidcc charles schwab investment management inc cut stock holding in bright horizon family solution inc nysebfam fluor corporation nyseflr share sold by vanguard group inc charles schwab investment management inc sell share of mattel inc nasdaqmat clearwater analytics holding inc nysecwan share acquired by vanguard group inc agco corporation nyseagco share sold by vanguard group inc vanguard group inc reduces position in ollies bargain outlet holding inc nasdaqolli vanguard group inc sell share of first american financial corporation nysefaf charles schwab investment management inc reduces holding in ally financial inc nyseally u bancorp de acquires share of graphic packaging holding company nysegpk vanguard group inc reduces position in lantheus holding inc nasdaqlnth charles schwab investment management inc boost position in dropbox inc nasdaqdbx vanguard group inc cut position in lincoln national corporation nyselnc vanguard group inc decrease position in itron inc nasdaqitri vanguard group inc trim stake in dropbox inc nasdaqdbx charles schwab investment management inc sell share of spotify technology nysespot chart industry inc nysegtls share purchased by charles schwab investment management inc nov inc nysenov share bought by vanguard group inc charles schwab investment management inc ha million stock position in spx technology inc nysespxc charles schwab investment management inc reduces position in talen energy corporation nasdaqtln curtisswright corporation nysecw position increased by charles schwab investment management inc vanguard group inc sell share of mattel inc nasdaqmat u bancorp de reduces holding in bce inc nysebce vanguard group inc lower stock holding in chart industry inc nysegtls atrium wealth solution inc ha stock position in ppg industry inc nyseppg analyst issue forecast for saputos q earnings tsesap atrium wealth solution inc buy share of vanguard international dividend appreciation etf nasdaqvigi vanguard group inc sell share of antero midstream corporation nyseam atrium wealth solution inc ha stake in southwest gas corporation nyseswx atrium wealth solution inc reduces position in sarepta therapeutic inc nasdaqsrpt vanguard group inc ha million holding in simpson manufacturing company inc nysessd vanguard group inc sell share of universal display corporation nasdaqoled vanguard group inc decrease stock holding in flowserve corporation nysefls share in spdr sp capital market etf nysearcakce bought by atrium wealth solution inc vanguard group inc sell share of credo technology group holding ltd nasdaqcrdo vanguard group inc acquires share of american airline group inc nasdaqaal atrium wealth solution inc invests in ugi corporation nyseugi atrium wealth solution inc decrease position in ingersoll rand inc nyseir home bancshares inc nysehomb share sold by atrium wealth solution inc vanguard group inc sell share of option care health inc nasdaqopch atrium wealth solution inc ha stock holding in rb global inc nyserba sl green realty corporation nyseslg share bought by vanguard group inc atrium wealth solution inc make new investment in crown holding inc nysecck notice regarding uaborkela financial statement for month period ended june nagarro release unaudited q result post yoy revenue growth in constant currency and increase in gross profit despite global macroeconomic challenge parent dont want cell phone ban at school we want smarter rule opinion mexico under pressure from trump sends cartel member to u the key to success at trumpputin alaska summit on ukraine low expectation opinion welcome to the menodivorce why woman arent sweating marriage in a sea of hot flash a tropical storm erin track west palm beach advises resident to stay alert prepared letter new challenge require a new congressman cctv foreign journalist open blind box at shuiting gate america must secure social security now republican pitch trump domestic policy agenda in iowa but some entrepreneur arent yet sold cnn politics the latest chatgpt is supposed to be phd level smart it cant even label a map wfp is leveraging it power in the mayoral race firstofitskind micropathway initiative now help community college build industryaligned training program my legacy is of broken men michael thomas on dream alcoholism and black fatherhood conquistador stolen yearold receipt return to mexico priscilla presley accused of hastening lisa mary death for financial gain memoir on a hill finding the best way to tell the story of america letter why now and why all the surprise eye on art frazier brings lowell to life in all gallery show rosner appointed ferc chair a trump move to reshape regulator u deploys force to take on latin american drug cartel in the southern caribbean today in history august fdr sign social security act now jimmy fallon is using taylor swift to suggest trump had sinister tie to epstein activist amariyanna copenys hometown crossword clue it not just dc republican seem happy to let trump do whatever he want opinion letter wondering if congressman still walk with jesus cincinnati beating suspect face year in prison after new charge from doj no peace without ukraine trumpputin alaska talk call for sober strategy oil city civic center to offer tour of national transit building thompson cosponsors bill to expand k career counseling idea opinion tom durkin come together debence receives grant for kid program grove city college ranked no for most religious student in u year ago aug lucas foley go to bat for bay state ag should too brokerage set cargojet inc tsecjt pt at c brokerage set penn entertainment inc nasdaqpenn target price at cbre group inc nysecbre receives consensus target price from analyst under armour nyseuaa price target lowered to at telsey advisory group john hood german settler helped shape nc rebellion today in historyaug failed gop candidate sentenced to year in connection with shooting spree at democrat official home stage mechanic and hanover in worcester make going out fun prokidney corp nasdaqprok receives average rating of hold from analyst spx technology inc nysespxc given average rating of buy by brokerage adapthealth corp nasdaqahco receives consensus target price from analyst contineum therapeutic inc nasdaqctnm receives consensus price target from brokerage on inc nyseontf receives consensus rating of hold from brokerage ap news summary at am edt douglas cohn and eleanor clift the shamed the shameless and the shirked editorial why is healey hesitating to call for tompkins resignation guest opinion erich menzel sculpture and art in the park event put loveland on the map ap news summary at am edt morning brief arrest since trump dc takeover trumpputin peace talk immigrant exodus social security is the cornerstone of retirement income fbi agent are again pulled from their day job to address a trump priority cnn politics cost of living in california still a discussion at county level propalestinian protester hold demonstration outside ivp reader sound off on turning against israel misusing troop and reckless bikers sbdc guide resident on starting legal homebased food business icf international nasdaqicfi price target raised to hang feng technology innovation co ltd fofo announces august th ipo alaska a source of russian imperial nostalgia powerful sister of north korean leader denies removal of frontline speaker cal thomas washington dc then and now emera tseema price target raised to c expensify nasdaqexfy cut to hold at wall street zen emera tseema price target raised to c rosenblatt security cut hut nasdaqhut price target to irs reportedly confirms stimulus check for eligible american rumor of august payout debunked plan late summer rollout american healthcare reit nyseahr price target raised to trump eas space regulation in boost for elon musk spacex jeff bezos blue origin and more pearson nysepso downgraded to hold rating by wall street zen first interstate bancsystem nasdaqfibk cut to sell at wall street zen morgan stanley cut primo brand nyseprmb price target to speaking out no longer interested in preventing pain and suffering north korea denies removing border loudspeaker ct real estate investment trust tsecrtun price target raised to c national bankshares cut boralex tseblx price target to c constellation software tsecsu price target raised to c share buyback program approved by plus lonplus board of director hong kong lashed with another torrential downpour redwire nyserdw price target lowered to at canaccord genuity group trisura group tsetsu price target raised to c gdi integrated facility service tsegdi downgraded by td security to hold leon furniture tselnf price target raised to c power co of canada tsepow price target raised to c at desjardins savaria tsesis price target raised to c saputo tsesap price target raised to c letter the right betrays every value it claim to hold youngkin order investigation into claim staff at virginia school arranged abortion without parental consent ap news summary at am edt flower food flo expected to announce earnings on friday kodiak gas service nysekgs to repurchase million in outstanding share trump nobel peace prize push gain steam president clinch backing from pakistan cambodia israel ahead of meeting with putin gop senator spill on the illegal measure he take to deal with dc crime fear softcat plc lonsct insider buy in stock jpmorgan global core real asset limited lonjara insider sell in stock fresnillo lonfres stock price expected to rise royal bank of canada analyst say asiamet resource lonars receives house stock rating from shore capital jpmorgan global core real asset limited lonjara insider sell in stock battalion oil post q revenue drop netease stock drop after chinese video game firm miss sale target aebi schmidt group building momentum after merging with the shyft group strong order backlog additional synergy upside and commitment to delever ask a real estate pro how can we stop complaining neighbor harassing behavior ask a real estate pro how can we stop complaining neighbor harassing behavior bessent see nvidia amd china agreement a a model for others nov declares regular quarterly dividend permanent magnet market worth billion by exclusive report by marketsandmarketstm clean energy tax credit market participant remain bullish the wall street journal rare earth rally after department of defense set a price floor aery technology report record q fy result cash flow positive m net income turnaround and strongest start to a fiscal year in company history envirogold global retains cantor fitzgerald a financial advisor p announces dual listing on nyse texas dividendpaying dow jones stock to buy in august u and australia support for nova estelle project realme pro review is it worth a buy shattuck lab report second quarter financial result and recent business highlight medium advisory air canada to hold medium availability with senior executive today aimia report second quarter result reiterates aebitda guidance and lower holdco cost target for the year advance auto part report second quarter result and reaffirms full year sale operating margin and free cash flow guidance calibercos sale drop percent nd century q earnings snapshot aquestive therapeutic announces pricing of million underwritten offering of common stock aquestive therapeutic announces pricing of million underwritten offering of common stock deere report third quarter net income of billion deere report third quarter net income of billion bcg report earnings ap technology summarybrief at am edt ap technology summarybrief at am edt today heloc home equity loan rate august bybit eu tap xion for inaugural launchpool in the eu opening regulated access for m user battalion oil corporation announces second quarter financial and operating result ford next model t moment started in california battalion oil corporation announces second quarter financial and operating result pro kapital council approved consolidated interim report for ii quarter and month of unaudited your guide to the best meme coin for exponential return this year pro kapital council approved consolidated interim report for ii quarter and month of unaudited cd rate today august take home up to u stock future lower a investor await inflation job data tool that have transformed business practice for entrepreneur fennec pharmaceutical report second quarter financial result and provides business update fennec pharmaceutical report second quarter financial result and provides business update kari lake defends voa cut in court after warning from capitol hill hta establishes advisory board want public input for planning highyield saving account rate today august rate are steady money market interest rate today august earn up to down is the worst over for lululemon stock the answer may surprise you lulu announces new credit agreement with white oak commercial finance llc mortgage rate today august rate hold steady spotting a hidden ice flight gold rush hit trumpera washington from price to the oval office xbocom exchange launch native xbo token and staking program ke holding inc to report second quarter financial result on august eastern time aquestive therapeutic announces m strategic funding agreement with rtw to support the potential launch of anaphylmtm epinephrine sublingual film this might be cadillacs craziest concept car ever motorcom big bank show there not just one recipe for success water recycle and reuse market worth billion by exclusive report by marketsandmarketstm ai modelasaservice is rising whats new in ai landscape radware and leading latvian isp tet sign managed security service provider agreement family target toprated school district this moving seasonbut at a price way roadside help harbor insurance see boost in revenue and customer satisfaction true random number generator market to hit usd billion by at a cagr of sn insider some state revenue is down but not the lottery vena receives sixth consecutive trustradius tech care award inception growth acquisition limited announces additional contribution to trust account to extend business combination period zelluna asa invitation to second quarter result webcast presentation aquestive therapeutic announces m strategic funding agreement with rtw to support the potential launch of anaphylmtm epinephrine sublingual film xunlei announces unaudited financial result for the second quarter ended june waterly secures m to power the future of water asset data sample management and more funding from burnt island venture and emerald will multiply the waterly platform reach ke holding inc to report second quarter financial result on august eastern time nd century group report second quarter financial result radware and leading latvian isp tet sign managed security service provider agreement how to enable p console sharing and offline play lulu announces new credit agreement with white oak commercial finance llc tungray technology appoints henry guo a chief financial officer over blizzard worker in irvine join union a gamingindustry labor movement expands xunlei announces unaudited financial result for the second quarter ended june the olympics give sponsor venue naming right for first time letter the city should consider alternative to flawed link plan better buy nike stock v lululemon stock near week low should you buy uipath stock beatendown undervalued growth stock down to buy hand over fist right now in august success of twin ownership deal will come down to how the money is spent the topperforming bank with more than b of asset crewscope raise over m to boost productivity for construction and industrial field crew disney stock buy or sell tyton partner release driving toward a degree career readiness and financial pressure at the forefront of student support best stock to buy novo nordisk stock v unitedhealth stock thing you need to know if you buy medtronic today ambow education sign white house pledge to expand access to aisupported education for america youth world share are mixed ahead of meeting between trump and putin ap news here are today mortgage refinance rate august rate decrease weak production data pours cold water on eurozones hope for industrial recovery charles schwab investment management inc lower stock holding in flowserve corporation nysefls rising rollouts of continuation vehicle signal maturity in private equity fiber reinforced polymer composite market size to surpass usd billion by driven by rising demand in construction automotive and aerospace sn insider china etown claim applied material poached staff misused proprietary chip tech openeden appoints bank of new york police in southern california find worth of labubus stolen from warehouse tech giant spotlight rising opportunity in western new york swalwells shameless flipflop from roleplaying carjacking with wife to shrugging off dc crime surge wisconsin request fema assistance after last weekend storm flooding why would a jew vote for mamdani rabbit in colorado are sprouting horn here why douthat maha like maga populism constrained by complicated marriage to big biz letter republican majority in congress ha sold out letter dl is right place for cruise ship terminal keene man died of multiple gunshot wound after traffic stop standoff trump urging for epstein distraction will lead to shocker at putin summit michael wolff letter christian nationalism richly hypocritical editorial grab hawaii feral cat bane by the scruff new york allowed pot shop to open too close to school now they might have to move alaska to take center stage at trumpputin summit alaska to take center stage at trumpputin summit sherrod brown give democrat jolt of enthusiasm in ohio sherrod brown give democrat jolt of enthusiasm in ohio doe unveils fresh billion funding to boost u critical mineral alaska to take center stage at trumpputin summit sherrod brown give democrat jolt of enthusiasm in ohio couple convicted for harboring illegal alien at texas bakery sheldon h jacobson made in america is alive well and misunderstood spotting a hidden ice flight editorial opposing mayor brandon johnson taxincrease idea is not a criminal act how trump and putin relationship ha evolved since they first met eight year ago cnn politics how trump and putin relationship ha evolved since they first met eight year ago willie wilson even pharaoh knew his limit when it came to tax editorial organ donor should be dead first thomas c bowen why ranked choice is the best cure for gerrymandering california democrat prepare redistricting effort to counter texas gop california democrat prepare redistricting effort to counter texas gop great american medium announces a belle christmas a new original movie starring ash tsai marshall williams and julia reilly family target toprated school district this moving seasonbut at a price chinese climber dy on k after reaching summit letter trump power grab threatens constitutional order letter the city should consider alternative to flawed link plan ice deportation airline avelo relies on bluestate subsidy will dem governor do anything about it trump effect dc police break from sanctuary city policy will cooperate with ice letter mob mentality at it worst ruin a meeting ichthyology is the branch of zoology that involves the study of which animal letter raucous rude and a cause for embarrasment foreign ngo say new israeli rule keep them from delivering gaza aid letter childish attendee ruin an important event father of israeli hostage delivers ringing message destroying hamas must be the main objective reader speak maybe former uconn woman basketball player can buy the connecticut sun kluth the u will regret throwing india under the bus a failed experiment in democracy letter democrat is sickened by town hall behavior israel announces west bank settlement that right group say could imperil palestinian state letter congressman show a lack of accountability zelenskyy meet with uk starmer a europe brace for trumpputin summit letter civility is important but rage is understandable israel announces a settlement that critic say will effectively sever the west bank in two sometimes you just cant put a price on a penny north state voice divided harrisburg house dems reject senate gop transit budget bill newsoms california redistricting push set up a standoff with republicanled opposition california republican are forming plan to oppose newsoms redistricting push cnn politics commentary welcome to american politics without norm commentary bad bunny is loud but puerto rican voter have spoken on statehood letter ugly hostility win out over open conversation ct electric rate to drop in september here why letter civility need to start at the top aid group call on israel to end weaponization of aid in gaza in gaza journalism is a death sentence reader speak trump need only to find answer to one question fit rover popular unconventional local health project awarded by doncaster council n korea accuses south of lying to the public opinion death tax and the nightmare that follows in ct letter town hall just a gigantic temper tantrum arkansas death row inmate sue to stop nitrogen gas execution asks to be next laura loomers attack on rfk jr the latest maga feud explained foton philippine launched a full commercial ev lineup adventure ready how the g tripros bumper cam take offroading to the next level she womens expo to offer health information saturday cisco ai infrastructure order top b why cryptocurrency okb ha more than doubled this week lyme disease is spiking this biotech nears a treatment and a buy point clean truck partnership take another blow with ftc deal in wild week for the beleaguered pact why share of coreweave are plunging this week unlocking the benefit of affordable insurance ace the test how to boost credit score tron flip cardano in market cap but this subm coin could outrun both dogecoin news today dogecoin eye breakout a whale accumulate m doge and golden cross form cold wallet presale hit m bittensor jump dogecoin volume soar best crypto to buy in citgo swing to m profit a bidding war intensifies for parent company israeli and palestinian bridge builder coming together to end the war aramco sign billion jafurah midstream deal with international consortium led by global infrastructure partner ethereum is outpacing cardano price but both are tasting crypto dust when compared to the gain from this presale ace the blue live in this mississippi delta town chicago fact of the day these analyst revise their forecast on coherent rosen leading investor counsel encourages cto realty growth inc investor to secure counsel before important deadline in security class action cto ctopa supreme court allows enforcement of mississippi social medium age verification law texture and franklinwh partner to deliver full control of home energy system ethereum solana lead crypto social medium discussion what about bitcoin doj sec charge two men in million water vending machine ponzi scheme cnbc recycle your electronics in beaufort county aug a car accident in smalltown tennessee lead to u charge against a major mexican drug operation eagle brown come together to gloriously throw down morning beer mississippi may require age verification parental consent for social medium supreme court say joby v archer two air taxi rival race toward but who really ha liftoff can ai cure cancer expert say the answer is complicated romania push nato to step up black sea mine task force eli lilly to raise uk price of mounjaro for outofpocket patient the wall street journal get ultimate allinone digital protection for off from mcafee advent vestment symbolism and tradition in liturgy these analyst increase their forecast on performance food group following strong q earnings how citi and jp morgan blockchain move influence payment philippine first floating solar farm set a national blueprint for clean energy best laptop for accounting in panorama insurance associate endorses california insurance market reform act call on industry to unite behind critical ballot initiative cleanspace continues rapid growth strengthening it concepttocommercialization model supreme court allows enforcement of mississippi social medium age verification law supreme court allows enforcement of mississippi social medium age verification law u dollar bullish structure contrast with softness in precious metal will dogecoin reach in august here the latest doge price forecast brinker international analyst boost their forecast following upbeat q result sonoma pharmaceutical valuation premium hinge on it revenue mix why pagseguro stock plummeted today money is rotating away from big tech stock it a good sign barrons confused trump forgets huge detail of big beautiful bill in cringe oval office moment prsany announces individual award recipient and under honoree ck hutchison profit plunge on oneoff loss from merger canaan report record q profit growth how to navigate the biggest pitfall of a growing business euna solution named government technology management team of the year in stevie award for technology excellence provectus q earnings snapshot auto insurance hike explained how to save without sacrificing coverage fleet electrification is getting easier by the day cisco to rally around here are top analyst forecast for thursday check out what whale are doing with ktos unpacking the latest option trading trend in sweetgreen spotlight on ionq analyzing the surge in option activity looking into domino pizza recent short interest what doe the market think about blackrock a look into mt bank inc price over earnings how is the market feeling about crown holding ozak ai potential may surpass bitcoins high financial trap keeping you stuck frustrated and how to escape past is prologue rising container volume portend falling rate lyno ai early bird presale exceeds k token sold a price increase approach pi coin news today pi coin dip to a investor confidence wane amid bitcoin surge blockdags m presale top trending crypto list despite cosmos positive outlook and k bnb prediction toll brother apartment living and harris realty company llc announce the opening of piper a new luxury apartment community in norwalk connecticut toll brother apartment living and harris realty company llc announce the opening of piper a new luxury apartment community in norwalk connecticut judge dismisses count against u rep cuellar of texas move bribery trial to next year average rate on a year mortgage drop to lowest level since october werner recall ladder due to faulty locking mechanism lm funding america swing to q profit taylor swift economic pull make orange the new green for big brand taylor swift economic pull make orange the new green for big brand connecticut credit union data breach affect verdek rank no on the inc list of america fastestgrowing private company why terawulf stock is skyrocketing today usa reject global carbon tax on global shipping small pacific island are hypocrite for their role in icj climate advisory say friend of science society walmart give employee extra saving on grocery with expanded discount dogecoin price outlook doge might hit however new gaming platform rblk offer bigger longterm upside ease protocol roll out government and enterprisegrade blockchain tool for realworld application fda alert over cookware that may leak lead into food product list todaycom kroger settle with c wholesale after failed billion albertsons merger oscar health face analyst question on path to profitability survey show parent will pay more for packed school lunch this year onfolio holding inc announces second quarter financial result and provides corporate update onfolio holding inc announces second quarter financial result and provides corporate update forget solana sol here are cheaper token to buy for the next gain in analyst predict a x surge for pepe coin but investor are loading up on the new xrp while they still can ibd live qa stock list for thursday aug baltimore face an energy crisis reader commentary best ultrawide camera smartphones to buy in toll brother announces opening of new model home at carrara estate in eagle idaho laura loomers attack on rfk jr the latest maga feud explained supreme court decline to block state law limiting kid use of social medium for now smiling maga man tear down democrat illegal trump hitler sign on highway overpass in iowa watch mayor adam float state bill to involuntarily hospitalize individual for public substance use amid reelection push the ad council and sugar launch ad council entertainment bipartisan solution for afghan ally spark hope israeli and palestinian bridge builder coming together to end the war the blue live in this mississippi delta town rosebud deputy and mhp arrest suspect accused of stealing truck after dodiks historic removal bosnia sovereignty crisis still smolder woman who dumped their manchild partner what wa the final straw live now california governor state l
'@

$ws.Range("A38").Value = 45884
$ws.Range("A38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B38").Value = 6449.7998046875
$ws.Range("C38").Value = 6481.33984375
$ws.Range("D38").Value = 6441.85009765625
$ws.Range("E38").Value = 6477.3798828125
$ws.Range("F38").Value = 4575060000
$ws.Range("G38").Value = -0.0028971350972291
$ws.Range("H38").Value = -1
$ws.Range("I38").Value = $headline
$ws.Range("J38").Value = $headline
$ws.Range("K38").Value = 0.151
$ws.Range("L38").Value = 0.08599999999999999
$ws.Range("M38").Value = 0.763
$ws.Range("N38").Value = 0.9999
$ws.Range("O38").Value = 112
$ws.Range("P38").Value = 112
$ws.Range("Q38").Value = 0
$ws.Range("R38").Value = 335
$ws.Range("S38").Value = 84
$ws.Range("T38").Value = 522
$ws.Range("U38").Value = 202
$ws.Range("V38").Value = 83
$ws.Range("W38").Value = 46
$ws.Range("X38").Value = 161
